$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename version label in row 7 (v45-11800 -> v5-11800)
$ws.Range("A7").Value = "v5-11800"

# Row 8: new benchmark entry
$ws.Cells.Item(8, 1).Value = "v6-8600"
$ws.Cells.Item(8, 2).Value = 16
$ws.Cells.Item(8, 3).Value = "Fixed movement code (pacman is getting stuck now)"
$ws.Cells.Item(8, 4).Value = 22
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 400
$ws.Cells.Item(8, 10).Value = 400
$ws.Cells.Item(8, 11).Value = 400
$ws.Cells.Item(8, 12).Value = 40
$ws.Cells.Item(8, 13).Value = 40
$ws.Cells.Item(8, 14).Value = 40
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 0
$ws.Cells.Item(8, 18).Value = 49
$ws.Cells.Item(8, 19).Value = 49
$ws.Cells.Item(8, 20).Value = 49

# Row 9: new benchmark entry
$ws.Cells.Item(9, 1).Value = "v7-9300"
$ws.Cells.Item(9, 2).Value = 64
$ws.Cells.Item(9, 3).Value = "Better ghost detection, ghost distance added to reward, removed pellet reward"
$ws.Cells.Item(9, 4).Value = 64
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 1752.19
$ws.Cells.Item(9, 10).Value = 1170
$ws.Cells.Item(9, 11).Value = 2610
$ws.Cells.Item(9, 12).Value = 162.09
$ws.Cells.Item(9, 13).Value = 109
$ws.Cells.Item(9, 14).Value = 226
$ws.Cells.Item(9, 15).Value = 0.27
$ws.Cells.Item(9, 16).Value = 0
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = 86.05
$ws.Cells.Item(9, 19).Value = 64.6
$ws.Cells.Item(9, 20).Value = 136.5

# Row 10: new benchmark entry (write B/C before A so new shared strings land
# in the same order as the source file)
$ws.Cells.Item(10, 2).Value = "55 (64)"
$ws.Cells.Item(10, 3).Value = "Removed learning when dying, accumulated ghost distance"
$ws.Cells.Item(10, 1).Value = "v8-2500"
$ws.Cells.Item(10, 4).Value = 21
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 110
$ws.Cells.Item(10, 10).Value = 110
$ws.Cells.Item(10, 11).Value = 110
$ws.Cells.Item(10, 12).Value = 11
$ws.Cells.Item(10, 13).Value = 11
$ws.Cells.Item(10, 14).Value = 11
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 0
$ws.Cells.Item(10, 18).Value = 11.2
$ws.Cells.Item(10, 19).Value = 11.2
$ws.Cells.Item(10, 20).Value = 11.2

# Update the remembered selection to match the author's last position
$ws.Range("G14").Select()
